$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9473684210526315
$ws.Range("C2").Value = 0.8947368421052632
$ws.Range("D2").Value = 0.9473684210526315
$ws.Range("E2").Value = 0.9473684210526315
$ws.Range("F2").Value = 0.9298245614035088
$ws.Range("G2").Value = 0.9298245614035088
$ws.Range("H2").Value = 0.9473684210526315
